# Actualización automática del inventario: agrega el nuevo producto
# "Kit Reparo de placa de impresoras Epson" como fila 76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").Value = "AVWLLD"
$ws.Range("B76").Value = "Kit Reparo de placa de impresoras Epson"
$ws.Range("C76").Value = "E09A88GA + Transistores C6144  y A2222"
$ws.Range("D76").Value = 0
$ws.Range("E76").Value = 100000
$ws.Range("F76").Value = 19
$ws.Range("G76").Value = 0
$ws.Range("H76").Formula = "=(E76-D76)*G76"
$ws.Range("I76").Formula = "=D76*F76"
$ws.Range("J76").Value = 0
